$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: merge C12/D12 into "\4,968", shift remaining labels left by one column ---
$ws.Range("C12").Value = "\4,968"
$ws.Range("D12").Value = "注文者"
$ws.Range("E12").Value = ":"
$ws.Range("F12").Value = "次材・購買課"
$ws.Range("G12").Value = "工藤"
$ws.Range("H12").ClearContents()

# --- Row 15: "Plate" -> "Pate", merge I15/J15 into "\1,200|", shift remaining cells left ---
$ws.Range("C15").Value = "Pate"
$ws.Range("I15").Value = "\1,200|"
$ws.Range("J15").Value = "9月10日"
$ws.Range("K15").Value = "|"
$ws.Range("L15").Value = "池田"
$ws.Range("M15").Value = "|9/5ヤグチ精機様持込み"
$ws.Range("N15").Value = "間"
$ws.Range("O15").ClearContents()

# --- Row 16: "Plate" -> "Pate", fix I16 value (no column shift needed) ---
$ws.Range("C16").Value = "Pate"
$ws.Range("I16").Value = "\1,200|"

# --- Row 17: "Plate" -> "Pate", merge I17/J17 into "\1,200|", shift remaining cells left ---
$ws.Range("C17").Value = "Pate"
$ws.Range("I17").Value = "\1,200|"
$ws.Range("J17").Value = "9月10日"
$ws.Range("K17").Value = "|"
$ws.Range("L17").Value = "池田"
$ws.Range("M17").Value = "|9/5ヤグチ精機様持込み"
$ws.Range("N17").Value = "拓"
$ws.Range("O17").ClearContents()

# --- Row 18: "Plate" -> "Pate" ---
$ws.Range("C18").Value = "Pate"
